$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 622.55554
$ws.Range("I6").Value = 391.33334
$ws.Range("K6").Value = 1174.00002
$ws.Range("M6").Value = -1062.00002

$ws.Range("H8").Value = 151.75
$ws.Range("I8").Value = 151.75
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 455.25
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -316.25
$ws.Range("N8").ClearContents()

$ws.Range("H17").Value = 1157.6666
$ws.Range("J17").Value = 1209.2
$ws.Range("L17").Value = 3627.6
$ws.Range("N17").Value = -3963.6

$ws.Range("H28").Value = 490.83334
$ws.Range("I28").Value = 490.83334
$ws.Range("K28").Value = 490.83334
$ws.Range("M28").Value = -5.833340000000021

$ws.Range("H52").Value = 750
$ws.Range("J52").Value = 750
$ws.Range("L52").Value = 2250
$ws.Range("N52").Value = -2570

$ws.Range("H62").Value = 9948
$ws.Range("I62").Value = 9948
$ws.Range("K62").Value = 9948
$ws.Range("M62").Value = -9324

$ws.Range("H65").Value = 9948
$ws.Range("I65").Value = 9948
$ws.Range("K65").Value = 49740
$ws.Range("M65").Value = -46620

$ws.Range("H80").Value = 517.3333
$ws.Range("J80").Value = 663.1667
$ws.Range("L80").Value = 1989.5001
$ws.Range("N80").Value = -3985.5001

$ws.Range("H83").Value = 517.3333
$ws.Range("J83").Value = 663.1667
$ws.Range("L83").Value = 5968.5003
$ws.Range("N83").Value = -15952.5003

$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -5812

$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -7808

$ws.Range("H107").Value = 605
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H118").Value = 575.7143
$ws.Range("I118").Value = 575.7143
$ws.Range("K118").Value = 1727.1429
$ws.Range("M118").Value = -70.14289999999983

$ws.Range("H125").Value = 2250
$ws.Range("J125").Value = 2250
$ws.Range("L125").Value = 20250
$ws.Range("N125").Value = -25170

$ws.Range("H138").Value = 10061.75
$ws.Range("J138").Value = 10061.75
$ws.Range("L138").Value = 30185.25
$ws.Range("N138").Value = -40465.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 9130.4
$ws.Range("I50").Value = 1848.6666
$ws.Range("J50").Value = 20053
$ws.Range("K50").Value = 1848.6666
$ws.Range("L50").Value = 20053
$ws.Range("M50").Value = -1134.6666
$ws.Range("N50").Value = -21481

$ws.Range("H102").Value = 166668560
$ws.Range("I102").Value = 166668560
$ws.Range("K102").Value = 166668560
$ws.Range("M102").Value = -166666938

$ws.Range("H122").Value = 41668332
$ws.Range("I122").Value = 41668332
$ws.Range("K122").Value = 125004996
$ws.Range("M122").Value = -125002546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 6325
$ws.Range("J36").Value = 2650
$ws.Range("L36").Value = 2650
$ws.Range("N36").Value = -3718

$ws.Range("H107").Value = 1193.5
$ws.Range("J107").Value = 1040.75
$ws.Range("L107").Value = 1040.75
$ws.Range("N107").Value = -4880.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 20.666666
$ws.Range("I7").Value = 9.666667
$ws.Range("J7").Value = 42.666668
$ws.Range("K7").Value = 9.666667
$ws.Range("L7").Value = 42.666668
$ws.Range("M7").Value = 103.333333
$ws.Range("N7").Value = -268.666668

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H31").Value = 5657.143
$ws.Range("I31").Value = 5602.75
$ws.Range("K31").Value = 5602.75
$ws.Range("M31").Value = -5307.75

$ws.Range("H34").Value = 5657.143
$ws.Range("I34").Value = 5602.75
$ws.Range("K34").Value = 5602.75
$ws.Range("M34").Value = -5400.75

$ws.Range("H41").Value = 12713.714
$ws.Range("J41").Value = 12713.714
$ws.Range("L41").Value = 12713.714
$ws.Range("N41").Value = -13569.714

$ws.Range("H95").Value = 12484.833
$ws.Range("J95").Value = 12484.833
$ws.Range("L95").Value = 12484.833
$ws.Range("N95").Value = -17976.833

$ws.Range("H122").Value = 1038.5
$ws.Range("I122").Value = 1038.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3115.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -665.5
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 5665.1816
$ws.Range("I132").Value = 5479.6665
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 16438.9995
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -13908.9995
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 101.5
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 9
$ws.Range("L46").Value = 600
$ws.Range("M46").Value = 82
$ws.Range("N46").Value = -782

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H108").Value = 568.7143
$ws.Range("I108").Value = 568.7143
$ws.Range("K108").Value = 1706.1429
$ws.Range("M108").Value = 1173.8571

$ws.Range("H109").Value = 2485.9092
$ws.Range("I109").Value = 1721.6666
$ws.Range("K109").Value = 5164.9998
$ws.Range("M109").Value = -4124.9998

$ws.Range("H117").Value = 3493.4
$ws.Range("J117").Value = 14250
$ws.Range("L117").Value = 42750
$ws.Range("N117").Value = -49634

$ws.Range("H131").Value = 2111
$ws.Range("I131").Value = 2499.5
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 7498.5
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = -2458.5
$ws.Range("N131").Value = -16080

$ws.Range("H132").Value = 2067.2
$ws.Range("I132").Value = 1120.6666
$ws.Range("J132").Value = 3487
$ws.Range("K132").Value = 10085.9994
$ws.Range("L132").Value = 31383
$ws.Range("M132").Value = -7555.999400000001
$ws.Range("N132").Value = -36443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1734.5714
$ws.Range("I113").Value = 1357
$ws.Range("K113").Value = 1357
$ws.Range("M113").Value = 813

$ws.Range("H126").Value = 1660
$ws.Range("I126").Value = 1490
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4470
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2000
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4715
$ws.Range("J40").Value = 5573.5
$ws.Range("L40").Value = 5573.5
$ws.Range("N40").Value = -5845.5

$ws.Range("H55").Value = 539.7143
$ws.Range("I55").Value = 475.8
$ws.Range("J55").Value = 699.5
$ws.Range("K55").Value = 475.8
$ws.Range("L55").Value = 699.5
$ws.Range("M55").Value = -302.8
$ws.Range("N55").Value = -1045.5

$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H122").Value = 55569652
$ws.Range("I122").Value = 66682184
$ws.Range("K122").Value = 200046552
$ws.Range("M122").Value = -200044102

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4749.5
$ws.Range("I3").Value = 4999
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 4999
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = -4885
$ws.Range("N3").Value = -4728

$ws.Range("H6").Value = 2051
$ws.Range("I6").Value = 1500
$ws.Range("J6").Value = 2602
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 2602
$ws.Range("M6").Value = -1385
$ws.Range("N6").Value = -2832

$ws.Range("H100").Value = 566.6667
$ws.Range("J100").Value = 900
$ws.Range("L100").Value = 1800
$ws.Range("N100").Value = -2882

$ws.Range("H107").Value = 2233
$ws.Range("I107").Value = 799
$ws.Range("K107").Value = 2397
$ws.Range("M107").Value = -477

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
